# Generate Report for Handback
# Adds "Latest Target File" (F) and "Latest Handback File" (G) values + hyperlinks
# to the zh-cn and de-de sheets, updates the Status column text, and refreshes the
# "Latest Handback DateTime" values to reflect the handback.

$wb = $excel.ActiveWorkbook

function Get-LinkAddress($ws, $addr) {
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $addr) {
            return $h.Address
        }
    }
    return $null
}

function Style-AsHyperlink($rng) {
    $rng.Font.Name = "Calibri"
    $rng.Font.Size = 11
    $rng.Font.Underline = $true
    $rng.Font.Color = 15570276  # RGB(100,149,237) == #6495ED, same blue as the existing HyperLink style
}

$statusText = "Handed back: in sync with en-US"

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Status column (C) - everything is now handed back and in sync with en-US
    $ws.Range("C2").Value = $statusText
    $ws.Range("C3").Value = $statusText

    # Latest Handback DateTime (H) - refresh with the handback timestamp
    if ($sheetName -eq "zh-cn") {
        $ws.Range("H2").Value = "2016-03-18 22:38:38"
        $ws.Range("H3").Value = "2016-03-18 22:38:38"
    } else {
        $ws.Range("H2").Value = "2016-03-18 22:38:43"
        $ws.Range("H3").Value = "2016-03-18 22:38:43"
    }

    # Row 2 - Latest Target File (F2) mirrors the source file (A2);
    #         Latest Handback File (G2) mirrors the handoff file (D2)
    $f2Display = $ws.Range("A2").Text
    $f2Address = Get-LinkAddress $ws '$A$2'
    $ws.Hyperlinks.Add($ws.Range("F2"), $f2Address, "", "", $f2Display) | Out-Null
    Style-AsHyperlink $ws.Range("F2")

    $g2Display = $ws.Range("D2").Text
    $g2Address = Get-LinkAddress $ws '$D$2'
    $ws.Hyperlinks.Add($ws.Range("G2"), $g2Address, "", "", $g2Display) | Out-Null
    Style-AsHyperlink $ws.Range("G2")

    # Row 3 - Latest Target File (F3) mirrors the source file (A3);
    #         Latest Handback File (G3) mirrors the handoff file (D3)
    $f3Display = $ws.Range("A3").Text
    $f3Address = Get-LinkAddress $ws '$A$3'
    $ws.Hyperlinks.Add($ws.Range("F3"), $f3Address, "", "", $f3Display) | Out-Null
    Style-AsHyperlink $ws.Range("F3")

    $g3Display = $ws.Range("D3").Text
    $g3Address = Get-LinkAddress $ws '$D$3'
    $ws.Hyperlinks.Add($ws.Range("G3"), $g3Address, "", "", $g3Display) | Out-Null
    Style-AsHyperlink $ws.Range("G3")
}

Write-Host "Handback report generated."
